# "error solve ifrs list"
#
# The IFRS consolidated figures for STX (rows 2-9 of company_list) were
# republished: most per-company rows shrink from "consolidated total"
# magnitudes down to much smaller corrected figures, two line items
# (당기순이익(비지배) / col J and 자본총계(비지배) / col O) are dropped for
# the first five data rows because they no longer apply, a trailing set of
# dividend-related columns (AG:AI) collapses to 0 for rows with values, and
# the last three data rows (7-9) lose every metric column beyond A:C
# (they keep only the rank/category/name columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> ordered column/value pairs that get a corrected value.
$rowValues = [ordered]@{
    2 = [ordered]@{
        D = 16346;  E = 357;    F = 357;     G = 3619;    H = 3806
        I = 3806;   K = 12550;  L = 11613;   M = 937;     N = 937
        P = 748;    Q = -761;   R = 450;     S = 770;     T = 15
        U = -776;   V = 8910;   W = 2.18;    X = 23.28;   Y = -182.91
        Z = 29.3;   AA = 1239.93; AB = 27.36; AC = 155647; AD = 0.22
        AE = 31403; AF = 1.07;  AG = 0;      AH = 0;      AI = 0
        AJ = 2984526
    }
    3 = [ordered]@{
        D = 16308;  E = -342;   F = -342;    G = -583;    H = -490
        I = -490;   K = 12672;  L = 11949;   M = 723;     N = 723
        P = 748;    Q = -160;   R = 421;     S = -618;    T = 10
        U = -170;   V = 9017;   W = -2.1;    X = -3;      Y = -59.01
        Z = -3.88;  AA = 1651.71; AB = -37.29; AC = -16405; AD = -2.09
        AE = 24227; AF = 1.42;  AG = 0;      AH = 0;      AI = 0
        AJ = 2987971
    }
    4 = [ordered]@{
        D = 17202;  E = -342;   F = -342;    G = -4579;   H = -4574
        I = -4574;  K = 7077;   L = 10383;   M = -3306;   N = -3306
        P = 1185;   Q = 578;    R = 532;     S = -1409;   T = 5
        U = 573;    V = 7910;   W = -1.99;   X = -26.59;  Y = 354.26
        Z = -46.33; AA = -314.07; AB = -377.22; AC = -145658; AD = -0.11
        AE = -69875; AF = -0.24; AG = 0;     AH = 0;      AI = 0
        AJ = 4733167
    }
    5 = [ordered]@{
        D = 18039;  E = 441;    F = 441;     G = 3346;    H = 3344
        I = 3344;   K = 6275;   L = 5850;    M = 426;     N = 426
        P = 478;    Q = 1003;   R = 473;     S = -962;    T = 4
        U = 998;    V = 3485;   W = 2.44;    X = 18.54;   Y = -232.18
        Z = 50.09;  AA = 1374.49; AB = 641.96; AC = 20929; AD = 1.14
        AE = 2226;  AF = 10.72; AG = 0;      AH = 0;      AI = 0
        AJ = 17862430
    }
    6 = [ordered]@{
        D = 15030;  E = -148;   F = -148;    G = 98;      H = 91
        I = 91;     K = 5541;   L = 4900;    M = 640;     N = 638
        P = 493;    Q = -481;   R = 248;     S = -288;    T = 6
        U = -487;   V = 3315;   W = -0.98;   X = 0.61;    Y = 17.12
        Z = 1.54;   AA = 765.34; AB = 644.53; AC = 471;   AD = 24.95
        AE = 3238;  AF = 3.63;  AJ = 19726902
    }
}

foreach ($row in $rowValues.Keys) {
    foreach ($col in $rowValues[$row].Keys) {
        $ws.Range("$col$row").Value = $rowValues[$row][$col]
    }
}

# Columns that are dropped outright (no longer reported) for rows 2-5:
# J = 당기순이익(비지배), O = 자본총계(비지배).
foreach ($row in 2..5) {
    $ws.Range("J$row").ClearContents()
    $ws.Range("O$row").ClearContents()
}

# Row 6 also drops the cash-dividend columns (AG:AI) entirely.
$ws.Range("AG6:AI6").ClearContents()

# Rows 7-9 lose every reported metric (D:AJ), keeping only rank/category/name
# (columns A:C).
foreach ($row in 7..9) {
    $ws.Range("D${row}:AJ${row}").ClearContents()
}

Write-Output "edit applied"
